# Critical category exibition fix
# Adds a new changelog entry row documenting the addition of a varchar
# "ano" column on the cashflowcategories table, and marks it as not
# committed to the server yet ("não").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new row's values first.
$ws.Range("A16").Value = 41832
$ws.Range("B16").Value = "Adicionado campo varchar ano na tabela cashflowcategories"
$ws.Range("C16").Value = "não"

# Copy the date formatting/style from the previous date cell (A15) so the
# new date cell A16 reuses the existing date style instead of creating a
# brand new one.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 41832

# Restore the selection that was active when the author saved the file.
$ws.Range("E24").Select()
